# Applies the row-content permutation described by the diff: the data
# (everything except the row number itself) that used to live in one row
# now lives in a different row. Row 11-24 form one cycle-set, rows 75-78
# another. Nothing else in the sheet changes.
#
# Strategy: snapshot every involved row's full set of relevant cell values
# (by column letter) BEFORE writing anything, then write the permuted
# snapshots back into place. Using Value2 for both read and write:
#   - read of a cell that doesn't exist returns $null
#   - write of $null clears/removes the cell (matches ClearContents)
# so cells that should disappear from a row (e.g. M, K, L, N, AC) do so
# automatically, and cells that should newly appear are (re)created.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All columns that ever hold data in the affected rows.
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY")

# target row -> source row: the NEW content of the target row is the OLD
# content that currently sits in the source row.
$mapping = @{
    11 = 13
    12 = 11
    13 = 14
    14 = 15
    15 = 12
    16 = 17
    17 = 20
    18 = 21
    19 = 22
    20 = 23
    21 = 24
    22 = 18
    23 = 19
    24 = 16
    75 = 76
    76 = 77
    77 = 78
    78 = 75
}

# All rows that participate (each row is both a source and a target here).
$involvedRows = @(11,12,13,14,15,16,17,18,19,20,21,22,23,24,75,76,77,78)

# 1) Snapshot current ("before") values for every involved row/column.
$snapshot = @{}
foreach ($r in $involvedRows) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range($col + $r).Value2
    }
    $snapshot[$r] = $rowData
}

# 2) Write each target row using the snapshot captured from its source row.
foreach ($target in $involvedRows) {
    $source = $mapping[$target]
    $srcData = $snapshot[$source]
    foreach ($col in $cols) {
        $ws.Range($col + $target).Value2 = $srcData[$col]
    }
}
